$wb = $excel.ActiveWorkbook

# --- Measures sheet: remove the unnecessary soil-reinforcement/stability-
# screen "2025" and "2045 met stabiliteitsscherm" measures for this
# DStability test case, then tidy the remaining max_outward (H) values.
$ws = $wb.Worksheets.Item("Measures")

# Delete from bottom to top so row indices of not-yet-deleted rows stay valid.
$ws.Rows.Item(5).Delete() | Out-Null
$ws.Rows.Item(4).Delete() | Out-Null
$ws.Rows.Item(3).Delete() | Out-Null

# Remaining 4 data rows (2-5) all get max_outward = 4 (was 30).
$ws.Range("H2:H5").Value = 4

# --- UI state: selection/active sheet moved to Measures!H6 while working on it.
$wsGeneral = $wb.Worksheets.Item("General")
$wsGeneral.Range("G18").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("H6").Select() | Out-Null
